$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering date range) ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Weekly crime-stat numeric updates (rows 14-30) ---
# Row 14
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = -75
$ws.Range("J14").Value = 38
$ws.Range("K14").Value = 2.631578947368
$ws.Range("M14").Value = -31.578947368421
$ws.Range("N14").Value = -74.838709677419

# Row 15
$ws.Range("C15").Value = 6
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 66.666666666666
$ws.Range("I15").Value = 138
$ws.Range("J15").Value = 137
$ws.Range("K15").Value = 0.729927007299
$ws.Range("L15").Value = 1.470588235294
$ws.Range("M15").Value = 22.123893805309
$ws.Range("N15").Value = -62.803234501347

# Row 16
$ws.Range("C16").Value = 28
$ws.Range("D16").Value = 39
$ws.Range("E16").Value = -28.205128205128
$ws.Range("F16").Value = 132
$ws.Range("G16").Value = 176
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 1110
$ws.Range("J16").Value = 1246
$ws.Range("K16").Value = -10.914927768860
$ws.Range("L16").Value = 34.220072551390
$ws.Range("M16").Value = -37.745372966909
$ws.Range("N16").Value = -87.324426173347

# Row 17
$ws.Range("C17").Value = 64
$ws.Range("D17").Value = 77
$ws.Range("E17").Value = -16.883116883116
$ws.Range("F17").Value = 307
$ws.Range("G17").Value = 322
$ws.Range("H17").Value = -4.658385093167
$ws.Range("I17").Value = 2251
$ws.Range("J17").Value = 2203
$ws.Range("K17").Value = 2.178847026781
$ws.Range("L17").Value = 18.473684210526
$ws.Range("M17").Value = 45.319561007101
$ws.Range("N17").Value = -48.945339079156

# Row 18
$ws.Range("C18").Value = 32
$ws.Range("D18").Value = 45
$ws.Range("E18").Value = -28.888888888888
$ws.Range("F18").Value = 134
$ws.Range("G18").Value = 192
$ws.Range("H18").Value = -30.208333333333
$ws.Range("I18").Value = 1086
$ws.Range("J18").Value = 1325
$ws.Range("K18").Value = -18.037735849056
$ws.Range("L18").Value = 4.022988505747
$ws.Range("M18").Value = -48.870056497175
$ws.Range("N18").Value = -90.651631230093

# Row 19
$ws.Range("C19").Value = 124
$ws.Range("D19").Value = 185
$ws.Range("E19").Value = -32.972972972973
$ws.Range("F19").Value = 555
$ws.Range("G19").Value = 653
$ws.Range("H19").Value = -15.007656967840
$ws.Range("I19").Value = 4188
$ws.Range("J19").Value = 4587
$ws.Range("K19").Value = -8.698495748855
$ws.Range("L19").Value = 38.721430937396
$ws.Range("M19").Value = 22.241681260945
$ws.Range("N19").Value = -24.567723342939

# Row 20
$ws.Range("C20").Value = 33
$ws.Range("D20").Value = 27
$ws.Range("E20").Value = 22.222222222222
$ws.Range("F20").Value = 155
$ws.Range("G20").Value = 140
$ws.Range("H20").Value = 10.714285714285
$ws.Range("I20").Value = 1166
$ws.Range("J20").Value = 1154
$ws.Range("K20").Value = 1.039861351819
$ws.Range("L20").Value = 50.064350064350
$ws.Range("M20").Value = -7.460317460317
$ws.Range("N20").Value = -92.205361320943

# Row 21
$ws.Range("C21").Value = 287
$ws.Range("D21").Value = 378
$ws.Range("E21").Value = -24.074074074074
$ws.Range("F21").Value = 1304
$ws.Range("G21").Value = 1499
$ws.Range("H21").Value = -13.008672448298
$ws.Range("I21").Value = 9978
$ws.Range("J21").Value = 10690
$ws.Range("K21").Value = -6.660430308699
$ws.Range("L21").Value = 29.014740108611
$ws.Range("M21").Value = -3.238944918541
$ws.Range("N21").Value = -78.223483195111

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -60
$ws.Range("F22").Value = 15
$ws.Range("H22").Value = 15.384615384615
$ws.Range("I22").Value = 115
$ws.Range("J22").Value = 121
$ws.Range("K22").Value = -4.958677685950
$ws.Range("L22").Value = 27.777777777777
$ws.Range("M22").Value = -38.502673796791

# Row 23
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 11
$ws.Range("E23").Value = -18.181818181818
$ws.Range("F23").Value = 51
$ws.Range("G23").Value = 61
$ws.Range("H23").Value = -16.393442622950
$ws.Range("I23").Value = 338
$ws.Range("J23").Value = 349
$ws.Range("K23").Value = -3.151862464183
$ws.Range("L23").Value = 10.819672131147
$ws.Range("M23").Value = 59.433962264150

# Row 24
$ws.Range("C24").Value = 318
$ws.Range("D24").Value = 333
$ws.Range("E24").Value = -4.504504504504
$ws.Range("F24").Value = 1279
$ws.Range("G24").Value = 1371
$ws.Range("H24").Value = -6.710430342815
$ws.Range("I24").Value = 10203
$ws.Range("J24").Value = 10047
$ws.Range("K24").Value = 1.552702299193
$ws.Range("L24").Value = 41.178912411789
$ws.Range("M24").Value = 31.668602400309

# Row 25
$ws.Range("C25").Value = 115
$ws.Range("D25").Value = 95
$ws.Range("E25").Value = 21.052631578947
$ws.Range("F25").Value = 509
$ws.Range("G25").Value = 445
$ws.Range("H25").Value = 14.382022471910
$ws.Range("I25").Value = 3758
$ws.Range("J25").Value = 3567
$ws.Range("K25").Value = 5.354639753294
$ws.Range("L25").Value = 25.685618729097
$ws.Range("M25").Value = -14.004576659038

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 22.727272727272
$ws.Range("I26").Value = 200
$ws.Range("J26").Value = 215
$ws.Range("K26").Value = -6.976744186046
$ws.Range("L26").Value = -6.542056074766

# Row 27
$ws.Range("F27").Value = 61
$ws.Range("G27").Value = 62
$ws.Range("H27").Value = -1.612903225806
$ws.Range("I27").Value = 431
$ws.Range("J27").Value = 451
$ws.Range("K27").Value = -4.434589800443
$ws.Range("L27").Value = 13.123359580052

# Row 28
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 14
$ws.Range("H28").Value = -39.130434782608
$ws.Range("I28").Value = 100
$ws.Range("J28").Value = 139
$ws.Range("K28").Value = -28.057553956834
$ws.Range("L28").Value = -22.480620155038
$ws.Range("M28").Value = -42.196531791907
$ws.Range("N28").Value = -80.879541108986

# Row 29
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 13
$ws.Range("H29").Value = -31.578947368421
$ws.Range("I29").Value = 87
$ws.Range("J29").Value = 108
$ws.Range("K29").Value = -19.444444444444
$ws.Range("L29").Value = -23.008849557522
$ws.Range("M29").Value = -38.732394366197
$ws.Range("N29").Value = -80.837004405286

# Row 30
$ws.Range("C30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 12
$ws.Range("H30").Value = -58.333333333333
$ws.Range("J30").Value = 86
$ws.Range("K30").Value = -44.186046511627
$ws.Range("L30").Value = 14.285714285714

$ws.Application.CutCopyMode = $false
